# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals.
# This updates column G ("K") values for rows 2-21 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 3
    4  = 4
    5  = 5
    6  = 12
    7  = 4
    8  = 12
    9  = 2
    10 = 1
    11 = 0
    12 = 2
    13 = 5
    14 = 3
    15 = 0
    16 = 0
    17 = 2
    18 = 2
    19 = 3
    20 = 5
    21 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
